$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.090.68"
$ws.Range("E2").Value = "'  +4.94%  "
$ws.Range("D3").Value = "'3.449.83"
$ws.Range("E3").Value = "'  +4.36%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  -0.10%  "
$ws.Range("D5").Value = "'582.52"
$ws.Range("E5").Value = "'  +5.88%  "
$ws.Range("D6").Value = "'184.09"
$ws.Range("E6").Value = "'  +7.25%  "
$ws.Range("D7").Value = "'0.631"
$ws.Range("E7").Value = "'  +2.64%  "
$ws.Range("D8").Value = "'3.444.05"
$ws.Range("E8").Value = "'  +4.49%  "
$ws.Range("E9").Value = "'  -0.06%  "
$ws.Range("E10").Value = "'  +1.90%  "
$ws.Range("D11").Value = "'0.643"
$ws.Range("E11").Value = "'  +2.64%  "
$ws.Range("D12").Value = "'56.21"
$ws.Range("E12").Value = "'  +5.76%  "
$ws.Range("D13").Value = "'0.0000277"
$ws.Range("E13").Value = "'  -0.04%  "
$ws.Range("D14").Value = "'9.41"
$ws.Range("E14").Value = "'  +4.82%  "
$ws.Range("D15").Value = "'4.003.55"
$ws.Range("E15").Value = "'  +4.07%  "
$ws.Range("D16").Value = "'18.61"
$ws.Range("E16").Value = "'  +3.42%  "
$ws.Range("D17").Value = "'3.450.81"
$ws.Range("E17").Value = "'  +4.68%  "
$ws.Range("D18").Value = "'67.010.97"
$ws.Range("E18").Value = "'  +4.71%  "
$ws.Range("E19").Value = "'  +0.60%  "
$ws.Range("D20").Value = "'12.07"
$ws.Range("E20").Value = "'  +4.04%  "
$ws.Range("D21").Value = "'1.01"
$ws.Range("E21").Value = "'  +3.55%  "
$ws.Range("D22").Value = "'480.98"
$ws.Range("E22").Value = "'  +6.17%  "
$ws.Range("D23").Value = "'5.48"
$ws.Range("E23").Value = "'  +10.91%  "
$ws.Range("D24").Value = "'16.66"
$ws.Range("D25").Value = "'4.40"
$ws.Range("E25").Value = "'  +8.99%  "
$ws.Range("D26").Value = "'89.40"
$ws.Range("E26").Value = "'  +3.40%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "'2.93"
$ws.Range("E27").Value = "'  +3.25%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.96"
$ws.Range("E28").Value = "'  +3.07%  "
$ws.Range("D29").Value = "'9.12"
$ws.Range("E29").Value = "'  +6.82%  "
$ws.Range("D30").Value = "'31.28"
$ws.Range("E30").Value = "'  +2.35%  "
$ws.Range("D31").Value = "'7.12"
$ws.Range("E31").Value = "'  +9.67%  "
$ws.Range("D32").Value = "'64.17"
$ws.Range("E32").Value = "'  +5.30%  "
$ws.Range("D33").Value = "'11.69"
$ws.Range("E33").Value = "'  +3.12%  "
$ws.Range("D34").Value = "'587.49"
$ws.Range("E34").Value = "'  +4.02%  "
$ws.Range("D35").Value = "'0.111"
$ws.Range("E35").Value = "'  +5.23%  "
$ws.Range("E36").Value = "'  +5.58%  "
$ws.Range("E37").Value = "'  +0.00%  "
$ws.Range("D38").Value = "'3.55"
$ws.Range("E38").Value = "'  +1.44%  "
$ws.Range("D39").Value = "'36.39"
$ws.Range("E39").Value = "'  +4.13%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.383"
$ws.Range("E40").Value = "'  +5.64%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "'0.0₃0767"
$ws.Range("E41").Value = "'  +6.35%  "
$ws.Range("D42").Value = "'3.205.55"
$ws.Range("E42").Value = "'  +5.71%  "
$ws.Range("D43").Value = "'2.89"
$ws.Range("E43").Value = "'  +5.85%  "
$ws.Range("D44").Value = "'0.0427"
$ws.Range("E44").Value = "'  +4.25%  "
$ws.Range("D45").Value = "'2.52"
$ws.Range("E45").Value = "'  +4.17%  "
$ws.Range("D46").Value = "'2.74"
$ws.Range("E46").Value = "'  +21.87%  "
$ws.Range("E47").Value = "'  +1.46%  "
$ws.Range("E48").Value = "'  +2.03%  "
$ws.Range("D49").Value = "'8.73"
$ws.Range("E49").Value = "'  +8.22%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "'  -0.10%  "
$ws.Range("D51").Value = "'3.20"
$ws.Range("E51").Value = "'  +10.79%  "
